# Apply the updated crypto price/volume figures to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string into a cell while forcing it to stay plain text
# (avoids Excel auto-converting a numeric-looking string, e.g. "26.83", into a
# real number) and then resets the cell style back to the default so no
# spurious "quote-prefix"/text-format style index is left behind on the cell.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '72.705.52'
$ws.Range("E2").Value = '  +4.08%  '
$ws.Range("D3").Value = '2.634.36'
$ws.Range("E3").Value = '  +2.68%  '
$ws.Range("E4").Value = '  +0.04%  '
Set-TextValue "D5" '605.24'
Set-TextValue "D6" '179.11'
$ws.Range("E6").Value = '  +0.52%  '
$ws.Range("E8").Value = '  +1.62%  '
$ws.Range("E9").Value = '  +8.51%  '
$ws.Range("D10").Value = '2.633.31'
$ws.Range("E10").Value = '  +2.74%  '
$ws.Range("E11").Value = '  +1.28%  '
Set-TextValue "D12" '0.355'
$ws.Range("E12").Value = '  +3.28%  '
Set-TextValue "D13" '5.03'
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("E14").Value = '  +4.37%  '
$ws.Range("E15").Value = '  +3.41%  '
$ws.Range("D16").Value = '72.472.85'
$ws.Range("E16").Value = '  +3.86%  '
Set-TextValue "D17" '26.83'
$ws.Range("E17").Value = '  +2.06%  '
$ws.Range("D18").Value = '2.628.64'
$ws.Range("E18").Value = '  +1.62%  '
Set-TextValue "D19" '11.78'
$ws.Range("E19").Value = '  +5.24%  '
Set-TextValue "D20" '385.32'
$ws.Range("E20").Value = '  +5.40%  '
Set-TextValue "D21" '7.92'
$ws.Range("E21").Value = '  +2.93%  '
$ws.Range("E22").Value = '  +1.51%  '
$ws.Range("E23").Value = '  +15.69%  '
Set-TextValue "D24" '74.16'
$ws.Range("E24").Value = '  +4.75%  '
$ws.Range("E25").Value = '  +3.06%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  +8.55%  '
$ws.Range("D28").Value = '2.730.20'
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  +4.23%  '
Set-TextValue "D31" '520.80'
$ws.Range("E31").Value = '  +0.92%  '
Set-TextValue "D32" '8.11'
$ws.Range("E32").Value = '  +4.09%  '
$ws.Range("E33").Value = '  +4.00%  '
$ws.Range("E34").Value = '  +1.55%  '
$ws.Range("E35").Value = '  -0.05%  '
Set-TextValue "D36" '163.02'
$ws.Range("E36").Value = '  -0.04%  '
Set-TextValue "D37" '19.41'
$ws.Range("E37").Value = '  +2.06%  '
$ws.Range("E38").Value = '  +3.97%  '
$ws.Range("E40").Value = '  -5.80%  '
$ws.Range("E41").Value = '  +5.42%  '
Set-TextValue "D42" '5.16'
$ws.Range("E42").Value = '  +4.58%  '
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("E44").Value = '  +4.59%  '
$ws.Range("E45").Value = '  +2.79%  '
Set-TextValue "D46" '39.45'
$ws.Range("E46").Value = '  +1.11%  '
Set-TextValue "D47" '151.01'
$ws.Range("E47").Value = '  -0.72%  '
Set-TextValue "D48" '3.70'
$ws.Range("E48").Value = '  +1.90%  '
$ws.Range("E49").Value = '  +4.42%  '
$ws.Range("E50").Value = '  +5.11%  '
$ws.Range("E51").Value = '  +2.48%  '
